$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 260; existing rows 260-279 shift down to 261-280
$ws.Rows.Item(260).Insert()

# Populate the newly inserted row 260 with the new weekly record
$ws.Cells.Item(260, 1).Value = 10
$ws.Cells.Item(260, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(260, 3).Value = "La Araucanía"
$ws.Cells.Item(260, 4).Value = 45021
$ws.Cells.Item(260, 5).Value = 9
$ws.Cells.Item(260, 6).Value = 100112005
$ws.Cells.Item(260, 7).Value = "Puerro"
$ws.Cells.Item(260, 8).Value = "Azul de Maquehue"
$ws.Cells.Item(260, 9).Value = "Primera"
$ws.Cells.Item(260, 10).Value = 45
$ws.Cells.Item(260, 11).Value = 12000
$ws.Cells.Item(260, 12).Value = 12000
$ws.Cells.Item(260, 13).Value = 12000
$ws.Cells.Item(260, 14).Value = "$/docena de paquetes"
$ws.Cells.Item(260, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(260, 16).Value = 1000
$ws.Cells.Item(260, 17).Value = 12
$ws.Cells.Item(260, 18).Value = "Hortaliza"
